$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update LoginDetails (sheet1) rows: row 6 becomes admin/password (what used
# to be row 7), and row 7 is removed entirely. ---
$ws1.Range("A6").Value = "admin"
$ws1.Range("B6").Value = "password"
$ws1.Range("A7").ClearContents()
$ws1.Range("B7").ClearContents()

# --- Add new column C with error-message data. Order of first assignment
# controls the shared-string table ordering, so set C2 before C1. ---
$ws1.Range("C2").Value = "Password cannot be empty"
$ws1.Range("C1").Value = "ErrorMessage"
$ws1.Range("C3").Value = "Username cannot be empty"
$ws1.Range("C4").Value = "Username cannot be empty"
$ws1.Range("C5").Value = "Invalid credentials"
$ws1.Range("C6").Value = "Invalid credentials"

# --- Add a new worksheet "ValidLoginDetails" after the existing sheet. ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "ValidLoginDetails"
$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "Admin"
$ws2.Range("B2").Value = "Qedge123!@#"
[void]$ws2.Range("F22").Select()

# --- Restore LoginDetails as the active sheet/selection. ---
$ws1.Activate()
[void]$ws1.Range("C1").Select()
